# Remove the second data row (Equipment Number "CSNU7211854", vessel
# "NORTHERN JUVENILE", etc.) from the sheet. Excel shifts all the rows
# below it up by one, so the old row 6 (HESU4027089 / NAVARINO / ...)
# disappears and the shared-strings table is recompacted to drop the
# strings that are no longer referenced by any cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:2").Delete()
